$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-of date range) ---
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

# Donor cells used to copy number formats (keeps existing style indices; avoids
# the COM layer fabricating brand new style entries when NumberFormat is set directly).
$style15 = $ws.Range("J23")
$style16 = $ws.Range("L26")
$style14 = $ws.Range("D26")

# Row 14
$style14.Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G14").Formula = '="0"'
$ws.Range("G14").Copy()
$ws.Range("G14").PasteSpecial(-4163)
$style14.Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("H14").Formula = '="***.*"'
$ws.Range("H14").Copy()
$ws.Range("H14").PasteSpecial(-4163)

# Row 15
$ws.Range("N15").Value = -93.333333333333

# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 500
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 150
$ws.Range("I16").Value = 35
$ws.Range("J16").Value = 13
$ws.Range("K16").Value = 169.230769230769
$ws.Range("L16").Value = 150
$ws.Range("M16").Value = 6.060606060606
$ws.Range("N16").Value = -78.395061728395

# Row 17
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = -3.846153846153
$ws.Range("I17").Value = 50
$ws.Range("J17").Value = 47
$ws.Range("K17").Value = 6.382978723404
$ws.Range("L17").Value = 38.888888888888
$ws.Range("M17").Value = 61.290322580645
$ws.Range("N17").Value = -46.808510638297

# Row 18
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 125
$ws.Range("I18").Value = 19
$ws.Range("J18").Value = 8
$ws.Range("K18").Value = 137.5
$ws.Range("L18").Value = 58.333333333333
$ws.Range("M18").Value = -26.923076923076
$ws.Range("N18").Value = -90.594059405940

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 21.428571428571
$ws.Range("I19").Value = 59
$ws.Range("J19").Value = 59
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 156.521739130435
$ws.Range("M19").Value = 34.090909090909
$ws.Range("N19").Value = -20.270270270270

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 150
$ws.Range("I20").Value = 19
$ws.Range("J20").Value = 7
$ws.Range("K20").Value = 171.428571428571
$ws.Range("L20").Value = 90
$ws.Range("M20").Value = -24
$ws.Range("N20").Value = -90.865384615384

# Row 21
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 25
$ws.Range("F21").Value = 98
$ws.Range("G21").Value = 71
$ws.Range("H21").Value = 38.028169014084
$ws.Range("I21").Value = 184
$ws.Range("J21").Value = 137
$ws.Range("K21").Value = 34.306569343065
$ws.Range("L21").Value = 87.755102040816
$ws.Range("M21").Value = 12.195121951219
$ws.Range("N21").Value = -75.693527080581

# Row 23
$ws.Range("C23").Value = 1
$style14.Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Formula = '="0"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$style14.Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Formula = '="***.*"'
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)
$ws.Range("F23").Value = 8
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 17
$ws.Range("K23").Value = 112.5
$ws.Range("L23").Value = 142.857142857143
$ws.Range("M23").Value = 750

# Row 24
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -19.047619047619
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 62
$ws.Range("H24").Value = 24.193548387096
$ws.Range("I24").Value = 161
$ws.Range("J24").Value = 111
$ws.Range("K24").Value = 45.045045045045
$ws.Range("L24").Value = 50.467289719626
$ws.Range("M24").Value = -9.039548022598

# Row 25
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -16.666666666666
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 56
$ws.Range("H25").Value = -1.785714285714
$ws.Range("I25").Value = 97
$ws.Range("J25").Value = 86
$ws.Range("K25").Value = 12.790697674418
$ws.Range("L25").Value = 110.869565217391
$ws.Range("M25").Value = -29.710144927536

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 10
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = 66.666666666666
$ws.Range("L27").Value = 66.666666666666

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 1
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -75
$ws.Range("N28").Value = -92.857142857142

# Row 29
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 1
$ws.Range("K29").Value = -50
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -75
$ws.Range("N29").Value = -91.666666666666

# Row 30
$ws.Range("L30").Value = -100

Write-Host "edits applied"
